# Kazakhstan Premier League - base update (19-06-2024 21:51)
# The edit re-orders a handful of match rows: the row metadata (column A,
# the sequential row index) stays put, but the rest of each row's data
# (columns B..AD: match id, teams, scores, odds, ...) gets shuffled between
# rows - effectively two row swaps and one 3-row rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

function Get-RowValues($ws, [int]$row, [int]$firstCol, [int]$lastCol) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($ws, [int]$row, [int]$firstCol, [int]$lastCol, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

# --- Swap rows 88 and 89 (data only, keep column A untouched) ---
$row88 = Get-RowValues $ws 88 $firstCol $lastCol
$row89 = Get-RowValues $ws 89 $firstCol $lastCol

Set-RowValues $ws 88 $firstCol $lastCol $row89
Set-RowValues $ws 89 $firstCol $lastCol $row88

# --- Rotate rows 98, 99, 100: 98 <- 99, 99 <- 100, 100 <- 98 ---
$row98  = Get-RowValues $ws 98  $firstCol $lastCol
$row99  = Get-RowValues $ws 99  $firstCol $lastCol
$row100 = Get-RowValues $ws 100 $firstCol $lastCol

Set-RowValues $ws 98  $firstCol $lastCol $row99
Set-RowValues $ws 99  $firstCol $lastCol $row100
Set-RowValues $ws 100 $firstCol $lastCol $row98

# --- Swap rows 156 and 157 (data only, keep column A untouched) ---
$row156 = Get-RowValues $ws 156 $firstCol $lastCol
$row157 = Get-RowValues $ws 157 $firstCol $lastCol

Set-RowValues $ws 156 $firstCol $lastCol $row157
Set-RowValues $ws 157 $firstCol $lastCol $row156
